$d = $word.ActiveDocument

$replacements = @(
    @{old="57×39="; new="70×75="},
    @{old="36×48="; new="41×59="},
    @{old="13×65="; new="61×20="},
    @{old="87×53="; new="69×90="},
    @{old="90×17="; new="71×38="},
    @{old="81×89="; new="95×58="},
    @{old="42×69="; new="11×83="},
    @{old="71×97="; new="88×73="},
    @{old="37×23="; new="46×36="},
    @{old="93×84="; new="93×19="},
    @{old="45×68="; new="78×21="},
    @{old="20×64="; new="23×69="},
    @{old="35×22="; new="51×60="},
    @{old="88×66="; new="59×40="},
    @{old="96×32="; new="88×97="},
    @{old="15×79="; new="62×92="},
    @{old="30×79="; new="85×71="},
    @{old="58×90="; new="26×63="},
    @{old="54×44="; new="51×71="},
    @{old="55×57="; new="28×25="},
    @{old="25×79="; new="20×54="},
    @{old="43×12="; new="91×90="},
    @{old="39×16="; new="16×53="},
    @{old="28×86="; new="22×89="},
    @{old="43×74="; new="93×88="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
